$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing "<end>" marker row (row 4) down to row 5 so a new
# Germination data row can be inserted in its place. Insert() preserves
# the moved row's original formatting.
$ws.Rows(4).Insert()

# Fill in the new data row 4 (3rd Germination sample)
$ws.Cells.Item(4, 1).Value = 3
$ws.Cells.Item(4, 2).Value = 2
$ws.Cells.Item(4, 3).Value = 400
$ws.Cells.Item(4, 4).Value = "18/12/2025"
$ws.Cells.Item(4, 5).Value = 13
$ws.Cells.Item(4, 6).Value = 2
$ws.Cells.Item(4, 7).Value = 1
$ws.Cells.Item(4, 8).Value = 3

# Re-apply the "date text" style (same as D2/D3) to D4, matching the rest
# of the Tgl Botol column
$ws.Range("D3").Copy()
$ws.Range("D4").PasteSpecial(-4122) # xlPasteFormats

# Widen the new column I to match column H
$ws.Columns.Item(9).ColumnWidth = $ws.Columns.Item(8).ColumnWidth

# New header "Type" in column I, formatted like the other headers
$ws.Range("I1").Value = "Type"
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122) # xlPasteFormats

# New "Type" values for each data row, formatted like the bold centered
# marker-row style
$ws.Range("I2").Value = "Direct"
$ws.Range("I3").Value = "Liquid"
$ws.Range("I4").Value = "Germination"
$ws.Range("A5").Copy()
$ws.Range("I2:I4").PasteSpecial(-4122) # xlPasteFormats

$ws.Application.CutCopyMode = $false

$ws.Range("D10").Select()
